# Update NMA and MA coefficient tables to use M instead of mu
# Renames the shared "mu_N" labels (column H) to "M_N" on every sheet,
# and updates the recorded active-cell selection on each sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # weibull
$ws2 = $wb.Worksheets.Item(2)   # gompertz
$ws3 = $wb.Worksheets.Item(3)   # fracpoly1
$ws4 = $wb.Worksheets.Item(4)   # fracpoly2

function Rename-MuColumn($ws, $lastRow) {
    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 8)   # column H
        $txt = $cell.Value2
        if ($txt -and $txt.ToString().StartsWith("mu_")) {
            $suffix = $txt.ToString().Substring(3)
            $cell.Value = "M_" + $suffix
        }
    }
}

Rename-MuColumn $ws1 7
Rename-MuColumn $ws2 7
Rename-MuColumn $ws3 10
Rename-MuColumn $ws4 10

# Update the stored selections on each sheet.
$ws2.Range("H2").Select() | Out-Null
$ws3.Range("H1").Select() | Out-Null
$ws4.Range("E34").Select() | Out-Null

# Restore sheet1 (weibull) as the active tab/sheet and set its selection.
$ws1.Activate() | Out-Null
$ws1.Range("H2").Select() | Out-Null
